# Remove the trailing "Ver no Jupiter..." / "© 2020 ..." footer paragraphs
# (and the blank paragraph that used to separate them from the requirement
# text above), leaving the single blank paragraph that precedes the final
# page-break paragraph untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$verRange = $d.Content
$verRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$verPara = $verRange.Paragraphs(1)

# The blank paragraph immediately preceding it is also removed.
$blankPara = $verPara.Previous()

# Locate the copyright/footer paragraph that follows.
$copyRange = $d.Content
$copyRange.Find.Execute("Contact: luizeleno@usp.br", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$copyPara = $copyRange.Paragraphs(1)

# Delete everything from the start of the blank paragraph through the end
# of the copyright paragraph (including its paragraph mark) in one shot.
$deleteRange = $d.Range($blankPara.Range.Start, $copyPara.Range.End)
$deleteRange.Delete()
